$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1149.591
$ws.Range("I41").Value = 896.5
$ws.Range("J41").Value = 1592.5
$ws.Range("K41").Value = 896.5
$ws.Range("L41").Value = 1592.5
$ws.Range("M41").Value = -456.5
$ws.Range("N41").Value = -2472.5
$ws.Range("H112").Value = 2593.3684
$ws.Range("J112").Value = 2923.4375
$ws.Range("L112").Value = 8770.3125
$ws.Range("N112").Value = -10986.3125
$ws.Range("H116").Value = 2159.5293
$ws.Range("I116").Value = 1601.875
$ws.Range("J116").Value = 2655.2222
$ws.Range("K116").Value = 1601.875
$ws.Range("L116").Value = 2655.2222
$ws.Range("M116").Value = 1840.125
$ws.Range("N116").Value = -9539.2222
$ws.Range("H121").Value = 1130.9375
$ws.Range("J121").Value = 1106.3334
$ws.Range("L121").Value = 3319.0002
$ws.Range("N121").Value = -6813.0002
$ws.Range("H129").Value = 801.6579
$ws.Range("I129").Value = 406.16666
$ws.Range("K129").Value = 1218.49998
$ws.Range("M129").Value = 3781.50002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 646.2143
$ws.Range("I2").Value = 494.09677
$ws.Range("K2").Value = 494.09677
$ws.Range("M2").Value = -381.09677
$ws.Range("H32").Value = 4647.34
$ws.Range("I32").Value = 3944.5684
$ws.Range("K32").Value = 3944.5684
$ws.Range("M32").Value = -3657.5684
$ws.Range("H45").Value = 1367.75
$ws.Range("I45").Value = 1372.4
$ws.Range("K45").Value = 1372.4
$ws.Range("M45").Value = -995.4000000000001
$ws.Range("H54").Value = 12000
$ws.Range("J54").Value = 12000
$ws.Range("L54").Value = 12000
$ws.Range("N54").Value = -13538
$ws.Range("H63").Value = 19233440
$ws.Range("I63").Value = 2302.8223
$ws.Range("J63").Value = 142862190
$ws.Range("K63").Value = 2302.8223
$ws.Range("L63").Value = 142862190
$ws.Range("M63").Value = -1616.8223
$ws.Range("N63").Value = -142863562
$ws.Range("H66").Value = 19233440
$ws.Range("I66").Value = 2302.8223
$ws.Range("J66").Value = 142862190
$ws.Range("K66").Value = 11514.1115
$ws.Range("L66").Value = 714310950
$ws.Range("M66").Value = -8082.111499999999
$ws.Range("N66").Value = -714317814
$ws.Range("H74").Value = 2650.44
$ws.Range("I74").Value = 1973.5714
$ws.Range("K74").Value = 1973.5714
$ws.Range("M74").Value = -1099.5714
$ws.Range("H77").Value = 2650.44
$ws.Range("I77").Value = 1973.5714
$ws.Range("K77").Value = 9867.857
$ws.Range("M77").Value = -5499.857
$ws.Range("H116").Value = 646.2143
$ws.Range("I116").Value = 494.09677
$ws.Range("K116").Value = 494.09677
$ws.Range("M116").Value = 1799.90323
$ws.Range("H124").Value = 16147.714
$ws.Range("J124").Value = 16147.714
$ws.Range("L124").Value = 16147.714
$ws.Range("N124").Value = -25967.714

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 646.2143
$ws.Range("I3").Value = 494.09677
$ws.Range("K3").Value = 494.09677
$ws.Range("M3").Value = -380.09677
$ws.Range("H134").Value = 1401
$ws.Range("I134").Value = 1271.5
$ws.Range("J134").Value = 1724.75
$ws.Range("K134").Value = 3814.5
$ws.Range("L134").Value = 5174.25
$ws.Range("M134").Value = -1279.5
$ws.Range("N134").Value = -10244.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1175.9678
$ws.Range("J31").Value = 2024.875
$ws.Range("L31").Value = 2024.875
$ws.Range("N31").Value = -2614.875
$ws.Range("H34").Value = 1175.9678
$ws.Range("J34").Value = 2024.875
$ws.Range("L34").Value = 2024.875
$ws.Range("N34").Value = -2428.875
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H132").Value = 1573.8182
$ws.Range("I132").Value = 1097.3914
$ws.Range("K132").Value = 3292.1742
$ws.Range("M132").Value = -762.1741999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 4722.1113
$ws.Range("J74").Value = 4722.1113
$ws.Range("L74").Value = 14166.3339
$ws.Range("N74").Value = -16288.3339
$ws.Range("H77").Value = 4722.1113
$ws.Range("J77").Value = 4722.1113
$ws.Range("L77").Value = 42499.00169999999
$ws.Range("N77").Value = -53107.00169999999
$ws.Range("H107").Value = 4482.08
$ws.Range("J107").Value = 5994.1113
$ws.Range("L107").Value = 17982.3339
$ws.Range("N107").Value = -21822.3339
$ws.Range("H131").Value = 711.42426
$ws.Range("I131").Value = 440.5
$ws.Range("J131").Value = 748.7931
$ws.Range("K131").Value = 1321.5
$ws.Range("L131").Value = 2246.3793
$ws.Range("M131").Value = 3718.5
$ws.Range("N131").Value = -12326.3793
$ws.Range("H139").Value = 1760.2
$ws.Range("J139").Value = 1887.05
$ws.Range("L139").Value = 5661.15
$ws.Range("N139").Value = -15941.15

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1720
$ws.Range("I41").Value = 1720
$ws.Range("K41").Value = 1720
$ws.Range("M41").Value = -1365
$ws.Range("H123").Value = 21665.2
$ws.Range("J123").Value = 21665.2
$ws.Range("L123").Value = 21665.2
$ws.Range("N123").Value = -26565.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 685.53845
$ws.Range("I22").Value = 623.44446
$ws.Range("J22").Value = 825.25
$ws.Range("K22").Value = 623.44446
$ws.Range("L22").Value = 825.25
$ws.Range("M22").Value = -328.44446
$ws.Range("N22").Value = -1415.25
$ws.Range("H27").Value = 685.53845
$ws.Range("I27").Value = 623.44446
$ws.Range("J27").Value = 825.25
$ws.Range("K27").Value = 623.44446
$ws.Range("M27").Value = -516.44446
$ws.Range("N27").Value = -1039.25
$ws.Range("H46").Value = 5272.727
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 6875
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 6875
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -7251

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2296.1667
$ws.Range("I132").Value = 2102.0645
$ws.Range("K132").Value = 6306.193499999999
$ws.Range("M132").Value = -3776.193499999999
$ws.Range("H136").Value = 1452.7826
$ws.Range("I136").Value = 1220.9333
$ws.Range("K136").Value = 3662.7999
$ws.Range("M136").Value = -1112.7999
